$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the ALU instruction count (B6), formulas in D6/F6 recalc automatically
$ws.Range("B6").Value = 816

# Update the Other instruction count (B10), formulas in D10/F10 recalc automatically
$ws.Range("B10").Value = 422

# Update the selected cell/range shown in the sheet view
$ws.Range("E23").Select()
